{"js": "// Remove the two duplicated \"featured image\" paragraphs that immediately\n// follow the \"Introducci\u00f3n\" and \"Presentaci\u00f3n del proyecto\" Heading2\n// paragraphs. Each of these paragraphs contains nothing but a single\n// centered inline picture that duplicates the image already shown at the\n// top of the document / right after the heading.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Collect, for every paragraph, whether it contains an inline picture and\n// its (trimmed) text, so we can find the heading paragraphs and the\n// duplicated image paragraphs that directly follow them.\nconst infos = paragraphs.items.map((p) => {\n  const inlinePics = p.inlinePictures;\n  inlinePics.load(\"items\");\n  return { paragraph: p, inlinePics };\n});\nawait context.sync();\n\nconst targets = [];\nfor (let i = 0; i < infos.length; i++) {\n  const text = infos[i].paragraph.text.trim();\n  if (text === \"Introducci\u00f3n\" || text === \"Presentaci\u00f3n del proyecto\") {\n    const next = infos[i + 1];\n    if (next && next.inlinePics.items.length > 0 && next.paragraph.text.trim() === \"\") {\n      targets.push(next.paragraph);\n    }\n  }\n}\n\n// Delete from the end so earlier indices stay valid while we iterate.\nfor (let i = targets.length - 1; i >= 0; i--) {\n  targets[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the two duplicated \"featured image\" paragraphs that immediately\n# follow the \"Introducci\u00f3n\" and \"Presentaci\u00f3n del proyecto\" headings. Each\n# of those paragraphs contains nothing but a single centered inline picture\n# that duplicates an image already shown elsewhere in the document.\n$d = $word.ActiveDocument\n\n$headings = @(\"Introducci\u00f3n\", \"Presentaci\u00f3n del proyecto\")\n\n# Snapshot the (flat, top-to-bottom) paragraph collection once; we will\n# delete from the end backwards so earlier indices remain valid.\n$paras = @($d.Paragraphs)\n\n$toDelete = @()\nfor ($i = 0; $i -lt $paras.Count; $i++) {\n    $text = $paras[$i].Range.Text.Trim()\n    if ($headings -contains $text) {\n        if ($i + 1 -lt $paras.Count) {\n            $next = $paras[$i + 1]\n            if ($next.Range.InlineShapes.Count -gt 0 -and $next.Range.Text.Trim() -eq \"\") {\n                $toDelete += $i + 1\n            }\n        }\n    }\n}\n\n$toDelete = $toDelete | Sort-Object -Descending\nforeach ($idx in $toDelete) {\n    $paras[$idx].Range.Delete()\n}\n"}
